# Update the "changes being made" bullet list: several description lines are
# reworded/expanded, and a brand-new bullet ("Forgot Password" security note)
# is inserted between the "Removing a Forgot Password button" heading and the
# "Fixing spelling mistakes in the Audit Trail" heading.

$d = $word.ActiveDocument

function Replace-ExactText([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    # Assign directly to Range.Text (rather than using Find's Replacement)
    # so Word's smart-quote / autocorrect typing substitutions do not kick
    # in and straight apostrophes are preserved verbatim.
    $rng.Text = $newText
}

Replace-ExactText "When a new user gets created they now get assigned to Requestor and are not automatically set to active. The user will now have to send a separate request to become active." "New users are assigned 'Requestor' by default and must submit a separate request to gain active status. Ensures user access control and compliance with role-based security policies."

Replace-ExactText "Adding ability to see the Lot location on the Standards Used tab on a Work Order record" "Enhancing the Calibration Standards tab by displaying Lot locations for calibration records in the Standards Used section of Work Orders, ensuring compliance with calibration data storage requirements, including instrument categories, calibration limits, and standard tracking."

Replace-ExactText "On the View All Equipment Update Requests screen, the QA Review field is now hidden" "Adjusting visibility settings in the Equipment Update Requests module by restricting access to the QA Review field, enforcing data access controls per system permissions."

Replace-ExactText "On a Audit Trail report in the reports section, the comment column will now have a note if the comment contains more than 3000+ characters" "Enhancing Audit Trail report functionality by implementing a flagging mechanism when comments exceed 3000 characters, ensuring accurate logging and traceability in compliance with record-keeping standards."

Replace-ExactText "On a Work Order report in the reports section, the location is visible and matches the location visible on the print preview version of the report " "Enhancing CMMS Cognos reporting functionality by ensuring Work Order reports accurately reflect location data, improving compliance with system-generated reporting standards and traceability."

Replace-ExactText "The Audit Trail which can be accessed on most record pages by right clicking to view every action that was made for a record, the column titles spelling mistakes are now fixed." "Maintaining system audit trail integrity by correcting event log formatting inconsistencies in audit records. This ensures compliance with system-generated audit trail storage, IT Change Control policies, and time-zone-based log retention as configured by the CMMS application server."

Replace-ExactText "On the User Setup screen users can now see a Last Login field" "Implementing a Last Login timestamp field in the User Setup screen to enhance user authentication tracking and facilitate compliance with login monitoring policies."

Replace-ExactText "The ECDB menu that was previously visible to all users is now hidden to only certain user roles" "Restricting ECDB menu visibility to specific user roles based on access permissions, enforcing compliance with role-based access controls."

# Insert the new "Forgot Password" security bullet right after the
# "Removing  a Forgot Password button" heading line, as a line break (not a
# new paragraph) followed by the new sentence, matching the rest of the
# document's <w:br/>-separated layout.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Removing  a Forgot Password button", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find text: Removing  a Forgot Password button"
}
$rng2.Collapse(0)
$rng2.InsertAfter([char]11 + "Strengthening authentication security by removing the Forgot Password recovery function to enforce LDAP authentication policies, ensuring application login strictly adheres to corporate IT security policy and requires standard network username and password credentials")
